$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4000
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872

$ws.Range("H76").Value = 5310
$ws.Range("I76").Value = 6020
$ws.Range("J76").Value = 4600
$ws.Range("K76").Value = 6020
$ws.Range("L76").Value = 4600
$ws.Range("M76").Value = -5705
$ws.Range("N76").Value = -5230

$ws.Range("H77").Value = 4000
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360

$ws.Range("H79").Value = 5310
$ws.Range("I79").Value = 6020
$ws.Range("J79").Value = 4600
$ws.Range("K79").Value = 6020
$ws.Range("L79").Value = 4600
$ws.Range("M79").Value = -4928
$ws.Range("N79").Value = -6784

$ws.Range("H107").Value = 2692.84
$ws.Range("J107").Value = 2420.3
$ws.Range("L107").Value = 2420.3
$ws.Range("N107").Value = -6260.3

$ws.Range("H112").Value = 2398.6316
$ws.Range("J112").Value = 2551.4119
$ws.Range("L112").Value = 7654.2357
$ws.Range("N112").Value = -9870.235700000001

$ws.Range("H138").Value = 1887.06
$ws.Range("J138").Value = 1973.4823
$ws.Range("L138").Value = 5920.4469
$ws.Range("N138").Value = -16200.4469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5553.321
$ws.Range("I32").Value = 5418.1763
$ws.Range("K32").Value = 5418.1763
$ws.Range("M32").Value = -5131.1763

$ws.Range("H74").Value = 3533.4546
$ws.Range("I74").Value = 2142.3333
$ws.Range("J74").Value = 5202.8
$ws.Range("K74").Value = 2142.3333
$ws.Range("L74").Value = 5202.8
$ws.Range("M74").Value = -1268.3333
$ws.Range("N74").Value = -6950.8

$ws.Range("H77").Value = 3533.4546
$ws.Range("I77").Value = 2142.3333
$ws.Range("J77").Value = 5202.8
$ws.Range("K77").Value = 10711.6665
$ws.Range("L77").Value = 26014
$ws.Range("M77").Value = -6343.666499999999
$ws.Range("N77").Value = -34750

$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3232.5
$ws.Range("I86").Value = 3543.4285
$ws.Range("J86").Value = 2797.2
$ws.Range("K86").Value = 3543.4285
$ws.Range("L86").Value = 2797.2
$ws.Range("M86").Value = -2420.4285
$ws.Range("N86").Value = -5043.2

$ws.Range("H89").Value = 3232.5
$ws.Range("I89").Value = 3543.4285
$ws.Range("J89").Value = 2797.2
$ws.Range("K89").Value = 17717.1425
$ws.Range("L89").Value = 13986
$ws.Range("M89").Value = -12101.1425
$ws.Range("N89").Value = -25218

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1463.4375
$ws.Range("I58").Value = 1250.85
$ws.Range("K58").Value = 1250.85
$ws.Range("M58").Value = -1047.85

$ws.Range("H136").Value = 1463.4375
$ws.Range("I136").Value = 1250.85
$ws.Range("K136").Value = 3752.55
$ws.Range("M136").Value = -1202.55

$ws.Range("H138").Value = 95495.55499999999
$ws.Range("J138").Value = 95495.55499999999
$ws.Range("L138").Value = 95495.55499999999
$ws.Range("N138").Value = -105775.555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 341.42856
$ws.Range("I50").Value = 80
$ws.Range("J50").Value = 537.5
$ws.Range("K50").Value = 240
$ws.Range("L50").Value = 1612.5
$ws.Range("M50").Value = 241
$ws.Range("N50").Value = -2574.5

$ws.Range("H53").Value = 341.42856
$ws.Range("I53").Value = 80
$ws.Range("J53").Value = 537.5
$ws.Range("K53").Value = 240
$ws.Range("L53").Value = 1612.5
$ws.Range("M53").Value = 241
$ws.Range("N53").Value = -2574.5

$ws.Range("H96").Value = 8320
$ws.Range("J96").Value = 8320
$ws.Range("L96").Value = 24960
$ws.Range("N96").Value = -29078

$ws.Range("H113").Value = 612.3158
$ws.Range("I113").Value = 526.3333
$ws.Range("J113").Value = 652
$ws.Range("K113").Value = 1578.9999
$ws.Range("L113").Value = 1956
$ws.Range("M113").Value = 591.0001
$ws.Range("N113").Value = -6296

$ws.Range("H122").Value = 1627.0526
$ws.Range("J122").Value = 1670.2222
$ws.Range("L122").Value = 15031.9998
$ws.Range("N122").Value = -19931.9998

$ws.Range("H131").Value = 12821404
$ws.Range("J131").Value = 926.5753
$ws.Range("L131").Value = 2779.7259
$ws.Range("N131").Value = -12859.7259

$ws.Range("H132").Value = 994.5
$ws.Range("I132").Value = 924.8182
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 8323.363800000001
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -5793.363800000001
$ws.Range("N132").Value = -16310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22503678
$ws.Range("I70").Value = 22731070
$ws.Range("J70").Value = 22225756
$ws.Range("K70").Value = 22731070
$ws.Range("L70").Value = 22225756
$ws.Range("M70").Value = -22730800
$ws.Range("N70").Value = -22226296

$ws.Range("H73").Value = 22503678
$ws.Range("I73").Value = 22731070
$ws.Range("J73").Value = 22225756
$ws.Range("K73").Value = 22731070
$ws.Range("L73").Value = 22225756
$ws.Range("M73").Value = -22730134
$ws.Range("N73").Value = -22227628

$ws.Range("H126").Value = 1798.2858
$ws.Range("I126").Value = 1567.6
$ws.Range("J126").Value = 2375
$ws.Range("K126").Value = 4702.799999999999
$ws.Range("L126").Value = 7125
$ws.Range("M126").Value = -2232.799999999999
$ws.Range("N126").Value = -12065

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1935.4286
$ws.Range("I7").Value = 1865.9166
$ws.Range("J7").Value = 2352.5
$ws.Range("K7").Value = 1865.9166
$ws.Range("L7").Value = 2352.5
$ws.Range("M7").Value = -1753.9166
$ws.Range("N7").Value = -2576.5

$ws.Range("H61").Value = 1249.4667
$ws.Range("I61").Value = 1141
$ws.Range("K61").Value = 1141
$ws.Range("M61").Value = -939

$ws.Range("H113").Value = 1249.4667
$ws.Range("I113").Value = 1141
$ws.Range("K113").Value = 1141
$ws.Range("M113").Value = 1029

$ws.Range("H126").Value = 1935.4286
$ws.Range("I126").Value = 1865.9166
$ws.Range("J126").Value = 2352.5
$ws.Range("K126").Value = 5597.7498
$ws.Range("L126").Value = 7057.5
$ws.Range("M126").Value = -3127.7498
$ws.Range("N126").Value = -11997.5

$ws.Range("H132").Value = 3316.2104
$ws.Range("I132").Value = 3301
$ws.Range("J132").Value = 3333.111
$ws.Range("K132").Value = 9903
$ws.Range("L132").Value = 9999.332999999999
$ws.Range("M132").Value = -7373
$ws.Range("N132").Value = -15059.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 37574.75
$ws.Range("J109").Value = 33319
$ws.Range("L109").Value = 33319
$ws.Range("N109").Value = -36093

$ws.Range("H113").Value = 404.72726
$ws.Range("I113").Value = 239.11111
$ws.Range("K113").Value = 717.3333299999999
$ws.Range("M113").Value = 1452.66667

$ws.Range("H132").Value = 1632.2142
$ws.Range("I132").Value = 1268.16
$ws.Range("K132").Value = 3804.48
$ws.Range("M132").Value = -1274.48
